$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 previously held "Leser" in A10 - remove it, new content starts at row 11
$ws.Range("A10").Value = $null

# --- Row 11: Server ---
$ws.Range("A11").Value = "Server"
$ws.Range("B11").Value = "sql6002.site4now.net"
$ws.Range("B11").Font.Name = "Arial"
$ws.Range("B11").Font.Size = 10
$ws.Range("B11").Font.Color = 3355443

# --- Row 12: Database ---
$ws.Range("A12").Value = "Database"
$ws.Range("B12").Value = "DB_A46E6D_admin"
$ws.Range("E12").Value = "DB_A46E6D_pms"
$ws.Range("H12").Value = "DB_A46E6D_sponsor"

# --- Row 13: Eier ---
$ws.Range("A13").Value = "Eier"
$ws.Range("B13").Value = "DB_A46E6D_Admin_admin"
$ws.Range("C13").Value = "Roma1995"
$ws.Range("E13").Value = "DB_A46E6D_pms_admin"
$ws.Range("F13").Value = "Roma1995"
$ws.Range("H13").Value = "DB_A46E6D_Sponsor_admin"
$ws.Range("I13").Value = "Roma1995"

# --- Row 14: Hovedbruker ---
$ws.Range("A14").Value = "Hovedbruker"
$ws.Range("B14").Value = "AdminMember"
$ws.Range("C14").Value = "Maine1953"
$ws.Range("E14").Value = "PMSMember"
$ws.Range("F14").Value = "Maine1953"
$ws.Range("H14").Value = "crew"
$ws.Range("I14").Value = "Maine1953"

# --- Row 16/17: Admin ODBC connection strings ---
$ws.Range("B16").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=sql6002.site4now.net;UID=m314alta;PWD=Maine1953;APP=Microsoft Office;DATABASE=DB_A46E6D_Admin;"
$ws.Range("B17").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=sql6002.site4now.net;UID=m314alta;PWD=Maine1953;APP=Microsoft Office;DATABASE=DB_A46E6D_AdminTest;"

# --- Row 18/19: PMS ODBC connection strings ---
$ws.Range("E18").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=sql6002.site4now.net;UID=m314alta;PWD=Maine1953;APP=Microsoft Office;DATABASE=DB_A46E6D_PMS;"
$ws.Range("E19").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=sql6002.site4now.net;UID=m314alta;PWD=Maine1953;APP=Microsoft Office;DATABASE=DB_A46E6D_PMSTest;"

# --- Row 20: Sponsor ODBC connection string ---
$ws.Range("H20").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=sql6002.site4now.net;UID=m314alta;PWD=Maine1953;APP=Microsoft Office;DATABASE=DB_A46E6D_Sponsor;"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.85546875
$ws.Columns.Item(2).ColumnWidth = 24.85546875
$ws.Columns.Item(5).ColumnWidth = 30.5703125

# --- View: scroll + selection ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E8").Select()
